# Updated symbol list on Sat Dec 24 22:43:06 UTC 2022 with GitHub Actions
#
# This script updates the "Price" column (D) values for a set of rows and
# two "Volume(1h)" column (E) label strings, matching the upstream data
# refresh. All D-column values are stored as plain text in the workbook
# (not numbers), so each assignment is written with a leading apostrophe to
# force Excel to keep them as text instead of silently re-typing them as
# numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2  = "244.65"
    3  = "21.90"
    4  = "5.394"
    7  = "0.8140"
    9  = "0.1440"
    10 = "0.07456"
    12 = "0.03049"
    13 = "0.09417"
    14 = "4.007"
    15 = "0.001588"
    16 = "0.04812"
    17 = "0.0005944"
    18 = "0.005421"
    19 = "0.004154"
    20 = "0.0009901"
    21 = "3.651"
    22 = "6.436"
    24 = "0.3243"
    26 = "0.00008506"
    27 = "0.0002902"
    40 = "0.03999"
    41 = "0.006434"
    42 = "0.1075"
    43 = "0.002722"
    44 = "0.006381"
    45 = "0.00005244"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# Volume(1h) label text updates
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
